$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = -21.57179999999998
$ws.Range("A6").Value = -22.68620000000001
$ws.Range("A7").Value = -19.42989999999999
$ws.Range("C7").Value = -12.1635
$ws.Range("C12").Value = -10.7145
$ws.Range("D13").Value = -8.536399999999997
$ws.Range("D14").Value = -7.915899999999997
$ws.Range("C15").Value = -14.63159999999999
$ws.Range("A16").Value = -21.70999999999999
$ws.Range("D16").Value = -9.112500000000008
$ws.Range("D19").Value = -8.580099999999996
$ws.Range("A20").Value = -19.28959999999999
$ws.Range("C20").Value = -11.9297
$ws.Range("C21").Value = -11.85760000000001
$ws.Range("C22").Value = -12.9343
$ws.Range("D22").Value = -8.072700000000003
$ws.Range("C23").Value = -11.86210000000001
$ws.Range("A28").Value = -21.79419999999999
$ws.Range("A29").Value = -21.32489999999997
$ws.Range("C29").Value = -11.36420000000001
$ws.Range("A32").Value = -21.2216
$ws.Range("C34").Value = -11.11060000000001
$ws.Range("D36").Value = -8.677299999999995
$ws.Range("A40").Value = -20.3099
$ws.Range("C42").Value = -12.0797
$ws.Range("C43").Value = -13.16589999999999
$ws.Range("C44").Value = -14.07209999999999
$ws.Range("C45").Value = -13.84619999999999
$ws.Range("A46").Value = -22.03470000000001
$ws.Range("C46").Value = -13.889
$ws.Range("D46").Value = -8.631299999999996
$ws.Range("C50").Value = -13.94389999999999
$ws.Range("D50").Value = -8.006599999999999
$ws.Range("A51").Value = -21.8745
$ws.Range("C51").Value = -11.6282
$ws.Range("A52").Value = -22.07509999999999
$ws.Range("A57").Value = -22.64620000000002
$ws.Range("A59").Value = -21.9299
$ws.Range("A62").Value = -22.02810000000002
$ws.Range("A66").Value = -21.62999999999999
$ws.Range("C66").Value = -10.9317
$ws.Range("C67").Value = -10.9824
$ws.Range("A73").Value = -20.73210000000001
$ws.Range("A74").Value = -21.95889999999999
$ws.Range("C79").Value = -11.37150000000002
$ws.Range("C84").Value = -14.07959999999999
$ws.Range("A92").Value = -21.69489999999999
$ws.Range("C92").Value = -11.1303
$ws.Range("D95").Value = -8.286799999999996
$ws.Range("C97").Value = -11.3672
$ws.Range("D97").Value = -8.187799999999996
$ws.Range("A100").Value = -22.1333
